$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-09-18 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-19 Thursday", 2)

$d.Content.Find.Execute("877÷4=219, 1", $true, $false, $false, $false, $false, $true, 1, $false, "698÷7=99, 5", 2)
$d.Content.Find.Execute("263÷9=29, 2", $true, $false, $false, $false, $false, $true, 1, $false, "163÷4=40, 3", 2)
$d.Content.Find.Execute("296÷9=32, 8", $true, $false, $false, $false, $false, $true, 1, $false, "480÷2=240, 0", 2)
$d.Content.Find.Execute("357÷4=89, 1", $true, $false, $false, $false, $false, $true, 1, $false, "144÷9=16, 0", 2)
$d.Content.Find.Execute("379÷8=47, 3", $true, $false, $false, $false, $false, $true, 1, $false, "804÷4=201, 0", 2)

$d.Content.Find.Execute("645÷7=92, 1", $true, $false, $false, $false, $false, $true, 1, $false, "112÷8=14, 0", 2)
$d.Content.Find.Execute("627÷4=156, 3", $true, $false, $false, $false, $false, $true, 1, $false, "765÷5=153, 0", 2)
$d.Content.Find.Execute("838÷4=209, 2", $true, $false, $false, $false, $false, $true, 1, $false, "872÷5=174, 2", 2)
$d.Content.Find.Execute("102÷4=25, 2", $true, $false, $false, $false, $false, $true, 1, $false, "565÷5=113, 0", 2)
$d.Content.Find.Execute("998÷9=110, 8", $true, $false, $false, $false, $false, $true, 1, $false, "983÷2=491, 1", 2)

$d.Content.Find.Execute("417÷6=69, 3", $true, $false, $false, $false, $false, $true, 1, $false, "656÷3=218, 2", 2)
$d.Content.Find.Execute("915÷2=457, 1", $true, $false, $false, $false, $false, $true, 1, $false, "416÷3=138, 2", 2)
$d.Content.Find.Execute("170÷9=18, 8", $true, $false, $false, $false, $false, $true, 1, $false, "308÷8=38, 4", 2)
$d.Content.Find.Execute("177÷9=19, 6", $true, $false, $false, $false, $false, $true, 1, $false, "674÷4=168, 2", 2)
$d.Content.Find.Execute("556÷6=92, 4", $true, $false, $false, $false, $false, $true, 1, $false, "885÷9=98, 3", 2)

$d.Content.Find.Execute("826÷4=206, 2", $true, $false, $false, $false, $false, $true, 1, $false, "615÷4=153, 3", 2)
$d.Content.Find.Execute("583÷7=83, 2", $true, $false, $false, $false, $false, $true, 1, $false, "924÷2=462, 0", 2)
$d.Content.Find.Execute("219÷5=43, 4", $true, $false, $false, $false, $false, $true, 1, $false, "427÷7=61, 0", 2)
$d.Content.Find.Execute("725÷6=120, 5", $true, $false, $false, $false, $false, $true, 1, $false, "923÷6=153, 5", 2)
$d.Content.Find.Execute("535÷3=178, 1", $true, $false, $false, $false, $false, $true, 1, $false, "536÷5=107, 1", 2)

$d.Content.Find.Execute("744÷4=186, 0", $true, $false, $false, $false, $false, $true, 1, $false, "645÷9=71, 6", 2)
$d.Content.Find.Execute("213÷3=71, 0", $true, $false, $false, $false, $false, $true, 1, $false, "955÷2=477, 1", 2)
$d.Content.Find.Execute("370÷5=74, 0", $true, $false, $false, $false, $false, $true, 1, $false, "511÷5=102, 1", 2)
$d.Content.Find.Execute("896÷5=179, 1", $true, $false, $false, $false, $false, $true, 1, $false, "154÷7=22, 0", 2)
$d.Content.Find.Execute("383÷4=95, 3", $true, $false, $false, $false, $false, $true, 1, $false, "202÷4=50, 2", 2)
